$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 2.35
$ws.Range("I3").Value = 3.4
$ws.Range("J3").Value = 3.25
$ws.Range("L3").Value = 4.33
$ws.Range("AB3").Value = 9.5
$ws.Range("AD3").Value = 23
$ws.Range("AE3").Value = 26
$ws.Range("AJ3").Value = 81
$ws.Range("AN3").Value = 13
$ws.Range("S4").Value = 3.75
$ws.Range("T4").Value = 1.28
$ws.Range("G5").Value = 2.05
$ws.Range("I5").Value = 4.33
$ws.Range("AB5").Value = 8
$ws.Range("AD5").Value = 19
$ws.Range("AL5").Value = 7.5
$ws.Range("AQ5").Value = 51
$ws.Range("G6").Value = 2
$ws.Range("I6").Value = 4.33
$ws.Range("J6").Value = 2.75
$ws.Range("AB6").Value = 7.5
$ws.Range("AI6").Value = 23
$ws.Range("AJ6").Value = 101
$ws.Range("AO6").Value = 51
$ws.Range("G7").Value = 6.25
$ws.Range("I7").Value = 1.6
$ws.Range("M7").Value = 1.07
$ws.Range("N7").Value = 9
$ws.Range("U7").Value = 4
$ws.Range("V7").Value = 1.22
$ws.Range("AC7").Value = 19
$ws.Range("G9").Value = 1.42
$ws.Range("I9").Value = 9
$ws.Range("M9").Value = 1.1
$ws.Range("N9").Value = 7
$ws.Range("Y9").Value = 2.63
$ws.Range("Z9").Value = 1.44
$ws.Range("AG9").Value = 7
$ws.Range("AI9").Value = 29
$ws.Range("AN9").Value = 29
$ws.Range("AP9").Value = 81
$ws.Range("G10").Value = 3
$ws.Range("I10").Value = 2.55
$ws.Range("J10").Value = 3.75
$ws.Range("L10").Value = 3.4
$ws.Range("Y10").Value = 2.1
$ws.Range("Z10").Value = 1.67
$ws.Range("AA10").Value = 7
$ws.Range("AC10").Value = 12
$ws.Range("AE10").Value = 29
$ws.Range("AG10").Value = 6.5
$ws.Range("AL10").Value = 6.5
$ws.Range("AM10").Value = 11
$ws.Range("AO10").Value = 26
$ws.Range("AP10").Value = 26
$ws.Range("H11").Value = 2.5
$ws.Range("I11").Value = 2.9
$ws.Range("J11").Value = 3.85
$ws.Range("K11").Value = 1.72
$ws.Range("L11").Value = 3.75
$ws.Range("M11").Value = 1.19
$ws.Range("N11").Value = 4.15
$ws.Range("O11").Value = 1.78
$ws.Range("P11").Value = 1.93
$ws.Range("Q11").Value = 3.25
$ws.Range("R11").Value = 1.3
$ws.Range("U11").Value = 6
$ws.Range("W11").Value = 1.75
$ws.Range("X11").Value = 1.98
$ws.Range("Y11").Value = 2.4
$ws.Range("AA11").Value = 5.7
$ws.Range("AC11").Value = 12.5
$ws.Range("AD11").Value = 40
$ws.Range("AE11").Value = 40
$ws.Range("AG11").Value = 4.15
$ws.Range("AH11").Value = 5.4
$ws.Range("AJ11").Value = 175
$ws.Range("AM11").Value = 12.5
$ws.Range("AP11").Value = 37
$ws.Range("H12").Value = 2.52
$ws.Range("I12").Value = 2.67
$ws.Range("J12").Value = 4
$ws.Range("K12").Value = 1.75
$ws.Range("L12").Value = 3.45
$ws.Range("M12").Value = 1.17
$ws.Range("N12").Value = 4.4
$ws.Range("O12").Value = 1.7
$ws.Range("P12").Value = 2.05
$ws.Range("Q12").Value = 3
$ws.Range("R12").Value = 1.34
$ws.Range("U12").Value = 5.5
$ws.Range("V12").Value = 1.11
$ws.Range("W12").Value = 1.65
$ws.Range("X12").Value = 2.1
$ws.Range("Y12").Value = 2.25
$ws.Range("Z12").Value = 1.57
$ws.Range("AA12").Value = 6.4
$ws.Range("AB12").Value = 14.5
$ws.Range("AC12").Value = 12.5
$ws.Range("AD12").Value = 50
$ws.Range("AF12").Value = 65
$ws.Range("AG12").Value = 4.4
$ws.Range("AI12").Value = 19
$ws.Range("AJ12").Value = 150
$ws.Range("AL12").Value = 5.7
$ws.Range("AN12").Value = 11
$ws.Range("AP12").Value = 32
$ws.Range("AQ12").Value = 55
$ws.Range("G13").Value = 2.6
$ws.Range("H13").Value = 2.7
$ws.Range("I13").Value = 3.1
$ws.Range("K13").Value = 1.87
$ws.Range("L13").Value = 3.7
$ws.Range("M13").Value = 1.13
$ws.Range("N13").Value = 5.1
$ws.Range("O13").Value = 1.55
$ws.Range("P13").Value = 2.3
$ws.Range("Q13").Value = 2.6
$ws.Range("R13").Value = 1.44
$ws.Range("U13").Value = 4.65
$ws.Range("V13").Value = 1.16
$ws.Range("W13").Value = 1.55
$ws.Range("X13").Value = 2.3
$ws.Range("Y13").Value = 2.05
$ws.Range("Z13").Value = 1.7
$ws.Range("AA13").Value = 6.1
$ws.Range("AB13").Value = 11.5
$ws.Range("AD13").Value = 30
$ws.Range("AE13").Value = 28
$ws.Range("AG13").Value = 5.1
$ws.Range("AH13").Value = 5.4
$ws.Range("AI13").Value = 16.5
$ws.Range("AL13").Value = 7.1
$ws.Range("AN13").Value = 11.25
$ws.Range("AO13").Value = 40
$ws.Range("AP13").Value = 32
$ws.Range("AQ13").Value = 50
$ws.Range("I14").Value = 8.5
$ws.Range("L14").Value = 7.7
$ws.Range("Q14").Value = 1.98
$ws.Range("R14").Value = 1.75
$ws.Range("X14").Value = 2.7
$ws.Range("Z14").Value = 1.53
$ws.Range("AM14").Value = 55
$ws.Range("G15").Value = 1.73
$ws.Range("I15").Value = 4.75
$ws.Range("U15").Value = 3.4
$ws.Range("V15").Value = 1.33
$ws.Range("AA15").Value = 7
$ws.Range("AF15").Value = 26
$ws.Range("AH15").Value = 6.5
$ws.Range("AI15").Value = 15
$ws.Range("AL15").Value = 13
$ws.Range("Y16").Value = 2.2
$ws.Range("Z16").Value = 1.62
$ws.Range("AF16").Value = 34
$ws.Range("M17").Value = 1.1
$ws.Range("N17").Value = 7
$ws.Range("O17").Value = 1.5
$ws.Range("P17").Value = 2.63
$ws.Range("AR17").Value = 1.93
$ws.Range("AS17").Value = 1.93
$ws.Range("G18").Value = 1.62
$ws.Range("J18").Value = 2.25
$ws.Range("AK18").Value = 501
$ws.Range("G20").Value = 2.2
$ws.Range("I20").Value = 3.25
$ws.Range("O21").Value = 1.36
$ws.Range("P21").Value = 3.2
$ws.Range("J24").Value = 2.38
$ws.Range("O24").Value = 1.29
$ws.Range("P24").Value = 3.5
$ws.Range("Q24").Value = 1.93
$ws.Range("R24").Value = 1.88
$ws.Range("U24").Value = 3.25
$ws.Range("V24").Value = 1.33
$ws.Range("AI24").Value = 15
$ws.Range("AK24").Value = 251
$ws.Range("AL24").Value = 13
$ws.Range("M26").Value = 1.06
$ws.Range("N26").Value = 10
